# NYPD 69th Precinct weekly CompStat report - new crime data collected.
# Updates the report header (volume/report-week text), widens column E to
# match the new data, and refreshes the Crime Complaints figures for rows
# 15-28 (Murder .. TOTAL) with the latest weekly/28-day/YTD/2-year counts
# and their computed percent changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: bump the bulletin's Volume/Number and the report week's
# date range.
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = "Volume 31   Number  22"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  5/27/2024  Through  6/2/2024"

# ---------------------------------------------------------------------
# Column E needs to be a bit wider for the new figures (matches column H's
# width, which is the closest width this engine's column-width quantizer
# can produce to the target 7.433768 characters).
# ---------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# ---------------------------------------------------------------------
# Some rows' "Week to Date" prior-year count (column C) drops to zero,
# which this sheet represents as the literal text "0" (shared style/string
# used elsewhere in the sheet, e.g. D28) instead of the numeric style used
# for nonzero counts. Use Copy so the destination picks up the exact same
# style + shared text "0" as the existing donor cell, rather than minting
# a new numeric-formatted style.
# ---------------------------------------------------------------------
$donorZero = $ws.Cells.Item(14, 3)

# --- Row 15 (Murder) ---
$donorZero.Copy($ws.Cells.Item(15, 3))
$ws.Cells.Item(15, 5).Value = -100
$ws.Cells.Item(15, 7).Value = 3
$ws.Cells.Item(15, 8).Value = -33.333333333333
$ws.Cells.Item(15, 9).Value = 5
$ws.Cells.Item(15, 10).Value = 5
$ws.Cells.Item(15, 12).Value = -50
$ws.Cells.Item(15, 14).Value = -50

# --- Row 16 (Rape) ---
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 5).Value = 66.666666666666
$ws.Cells.Item(16, 6).Value = 9
$ws.Cells.Item(16, 8).Value = -30.769230769230
$ws.Cells.Item(16, 9).Value = 35
$ws.Cells.Item(16, 10).Value = 60
$ws.Cells.Item(16, 11).Value = -41.666666666666
$ws.Cells.Item(16, 12).Value = -28.571428571428
$ws.Cells.Item(16, 13).Value = -51.388888888888
$ws.Cells.Item(16, 14).Value = -85.82995951417

# --- Row 17 (Robbery) ---
$ws.Cells.Item(17, 3).Value = 5
$ws.Cells.Item(17, 4).Value = 6
$ws.Cells.Item(17, 5).Value = -16.666666666666
$ws.Cells.Item(17, 6).Value = 23
$ws.Cells.Item(17, 7).Value = 13
$ws.Cells.Item(17, 8).Value = 76.923076923076
$ws.Cells.Item(17, 9).Value = 88
$ws.Cells.Item(17, 10).Value = 90
$ws.Cells.Item(17, 11).Value = -2.222222222222
$ws.Cells.Item(17, 12).Value = -16.981132075471
$ws.Cells.Item(17, 13).Value = 41.935483870967
$ws.Cells.Item(17, 14).Value = -5.376344086021

# --- Row 18 (Fel. Assault) ---
$donorZero.Copy($ws.Cells.Item(18, 3))
$ws.Cells.Item(18, 5).Value = -100
$ws.Cells.Item(18, 6).Value = 5
$ws.Cells.Item(18, 8).Value = -44.444444444444
$ws.Cells.Item(18, 9).Value = 19
$ws.Cells.Item(18, 10).Value = 28
$ws.Cells.Item(18, 11).Value = -32.142857142857
$ws.Cells.Item(18, 12).Value = -54.761904761904
$ws.Cells.Item(18, 13).Value = -81.904761904761
$ws.Cells.Item(18, 14).Value = -93.189964157706

# --- Row 19 (Burglary) ---
$ws.Cells.Item(19, 3).Value = 4
$ws.Cells.Item(19, 4).Value = 3
$ws.Cells.Item(19, 5).Value = 33.333333333333
$ws.Cells.Item(19, 6).Value = 20
$ws.Cells.Item(19, 7).Value = 11
$ws.Cells.Item(19, 8).Value = 81.818181818181
$ws.Cells.Item(19, 10).Value = 91
$ws.Cells.Item(19, 11).Value = 41.758241758241
$ws.Cells.Item(19, 12).Value = 31.632653061224
$ws.Cells.Item(19, 13).Value = 41.758241758241
$ws.Cells.Item(19, 14).Value = -12.837837837837

# --- Row 20 (Gr. Larceny) ---
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Value = 3
$ws.Cells.Item(20, 5).Value = -66.666666666666
$ws.Cells.Item(20, 6).Value = 10
$ws.Cells.Item(20, 7).Value = 10
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 58
$ws.Cells.Item(20, 10).Value = 65
$ws.Cells.Item(20, 11).Value = -10.769230769230
$ws.Cells.Item(20, 12).Value = -10.769230769230
$ws.Cells.Item(20, 13).Value = 34.883720930232
$ws.Cells.Item(20, 14).Value = -93.605292171995

# --- Row 21 (G.L.A. / TOTAL) ---
$ws.Cells.Item(21, 3).Value = 15
$ws.Cells.Item(21, 4).Value = 17
$ws.Cells.Item(21, 5).Value = -11.764705882352
$ws.Cells.Item(21, 6).Value = 69
$ws.Cells.Item(21, 8).Value = 15
$ws.Cells.Item(21, 9).Value = 336
$ws.Cells.Item(21, 10).Value = 343
$ws.Cells.Item(21, 11).Value = -2.040816326530
$ws.Cells.Item(21, 12).Value = -9.677419354838
$ws.Cells.Item(21, 13).Value = -11.811023622047
$ws.Cells.Item(21, 14).Value = -80.094786729857

# --- Row 23 (Transit) ---
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = 5
$ws.Cells.Item(23, 5).Value = -60
$ws.Cells.Item(23, 6).Value = 10
$ws.Cells.Item(23, 7).Value = 14
$ws.Cells.Item(23, 8).Value = -28.571428571428
$ws.Cells.Item(23, 9).Value = 48
$ws.Cells.Item(23, 10).Value = 53
$ws.Cells.Item(23, 11).Value = -9.433962264150
$ws.Cells.Item(23, 12).Value = 20
$ws.Cells.Item(23, 13).Value = 152.631578947368

# --- Row 24 (Housing) ---
$ws.Cells.Item(24, 3).Value = 11
$ws.Cells.Item(24, 4).Value = 12
$ws.Cells.Item(24, 5).Value = -8.333333333333
$ws.Cells.Item(24, 6).Value = 43
$ws.Cells.Item(24, 7).Value = 48
$ws.Cells.Item(24, 8).Value = -10.416666666666
$ws.Cells.Item(24, 9).Value = 221
$ws.Cells.Item(24, 10).Value = 263
$ws.Cells.Item(24, 11).Value = -15.969581749049
$ws.Cells.Item(24, 12).Value = -4.741379310344
$ws.Cells.Item(24, 13).Value = 6.763285024154

# --- Row 25 (Petit Larceny) ---
$ws.Cells.Item(25, 3).Value = 3
$ws.Cells.Item(25, 4).Value = 5
$ws.Cells.Item(25, 5).Value = -40
$ws.Cells.Item(25, 6).Value = 9
$ws.Cells.Item(25, 7).Value = 14
$ws.Cells.Item(25, 8).Value = -35.714285714285
$ws.Cells.Item(25, 9).Value = 45
$ws.Cells.Item(25, 10).Value = 70
$ws.Cells.Item(25, 11).Value = -35.714285714285
$ws.Cells.Item(25, 12).Value = -55

# --- Row 26 (Retail Theft) ---
$ws.Cells.Item(26, 3).Value = 10
$ws.Cells.Item(26, 4).Value = 9
$ws.Cells.Item(26, 5).Value = 11.111111111111
$ws.Cells.Item(26, 6).Value = 30
$ws.Cells.Item(26, 7).Value = 34
$ws.Cells.Item(26, 8).Value = -11.764705882352
$ws.Cells.Item(26, 9).Value = 141
$ws.Cells.Item(26, 10).Value = 156
$ws.Cells.Item(26, 11).Value = -9.615384615384
$ws.Cells.Item(26, 12).Value = -9.615384615384
$ws.Cells.Item(26, 13).Value = -24.193548387096

# --- Row 27 (Misd. Assault) ---
$donorZero.Copy($ws.Cells.Item(27, 3))
$ws.Cells.Item(27, 5).Value = -100
$ws.Cells.Item(27, 9).Value = 5
$ws.Cells.Item(27, 10).Value = 7
$ws.Cells.Item(27, 11).Value = -28.571428571428
$ws.Cells.Item(27, 12).Value = -68.75

# --- Row 28 (UCR Rape*) ---
$donorZero.Copy($ws.Cells.Item(28, 3))
$ws.Cells.Item(28, 6).Value = 2
$ws.Cells.Item(28, 7).Value = 2
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 12).Value = 7.692307692307
